# Sideways - Movie review.docx edit script
$d = $word.ActiveDocument

# 1. Cast line: merge "Cast: Paul " + "Giamatti (" + "Miles), Thomas Haden Church" runs (no text change)
$d.Content.Find.Execute("Cast: Paul Giamatti (Miles), Thomas Haden Church", $false, $false, $false, $false, $false, $true, 1, $false, "Cast: Paul Giamatti (Miles), Thomas Haden Church", 2) | Out-Null

# 2. Cast line: merge "(Jack), " + "Virginia" + " Madsen" runs (no text change)
$d.Content.Find.Execute("(Jack), Virginia Madsen", $false, $false, $false, $false, $false, $true, 1, $false, "(Jack), Virginia Madsen", 2) | Out-Null

# 3. Plot paragraph edits
$d.Content.Find.Execute("teacher(Miles), who is", $true, $false, $false, $false, $false, $true, 1, $false, "teacher, Miles (who is", 2) | Out-Null
$d.Content.Find.Execute("book published, going", $true, $false, $false, $false, $false, $true, 1, $false, "book published) going", 2) | Out-Null
$d.Content.Find.Execute("actor friend(Jack), to celebrate that the friend is getting", $true, $false, $false, $false, $false, $true, 1, $false, "actor friend Jack, to celebrate that Jack is getting", 2) | Out-Null
$d.Content.Find.Execute("goes back home.", $true, $false, $false, $false, $false, $true, 1, $false, "go back home.", 2) | Out-Null

# 4. "wastes" -> "waste"; "them, already" -> "them Already"
$d.Content.Find.Execute("Jack wastes no time", $true, $false, $false, $false, $false, $true, 1, $false, "Jack waste no time", 2) | Out-Null
$d.Content.Find.Execute("dislike them, already from", $true, $false, $false, $false, $false, $true, 1, $false, "dislike them Already from", 2) | Out-Null

# 5. Para 7: "Jack is almost the opposite of Miles..."
$d.Content.Find.Execute("Miles, he is really positive, we then meet", $true, $false, $false, $false, $false, $true, 1, $false, "Miles; he is really positive and optimistic. We then meet", 2) | Out-Null
$d.Content.Find.Execute("serving tasting wines,", $true, $false, $false, $false, $false, $true, 1, $false, "serving wine tastings,", 2) | Out-Null
$d.Content.Find.Execute("Miles and Stephanie is negative", $true, $false, $false, $false, $false, $true, 1, $false, "Miles and Stephanie are negative", 2) | Out-Null

# 6. Para 8: "as soon as the switch" -> "as soon as they switch"
$d.Content.Find.Execute("as soon as the switch to red wine", $true, $false, $false, $false, $false, $true, 1, $false, "as soon as they switch to red wine", 2) | Out-Null

# 7. Para 10: "finally begin to establish" -> "finally beginning to establish"
$d.Content.Find.Execute("finally begin to establish", $true, $false, $false, $false, $false, $true, 1, $false, "finally beginning to establish", 2) | Out-Null

# 8. Para 11: "great about having us invested" -> "great about getting us invested"
$d.Content.Find.Execute("great about having us invested", $true, $false, $false, $false, $false, $true, 1, $false, "great about getting us invested", 2) | Out-Null
# "his book not and gets really upset, we sympathize" -> "his book not being published and gets really upset, and we sympathize"
$d.Content.Find.Execute("his book not and gets really upset, we sympathize", $true, $false, $false, $false, $false, $true, 1, $false, "his book not being published and gets really upset, and we sympathize", 2) | Out-Null

# 9. Para 14: merge runs + "red/white theme" -> "red/white contrast" + remove proofErr around "don't"
$d.Content.Find.Execute("I really like the red/white theme in the movie, and how like in Miles book, we don’t know if he gets together with Maya or what happens to him. ", $true, $false, $false, $false, $false, $true, 1, $false, "I really like the red/white contrast in the movie, and how like in Miles book, we don’t know if he gets together with Maya or what happens to him. ", 2) | Out-Null

# 10. Para 12: merge runs (remove proofErr around "we're"), no text change
$d.Content.Find.Execute("and the next day we’re back to the black background", $true, $false, $false, $false, $false, $true, 1, $false, "and the next day we’re back to the black background", 2) | Out-Null

# 11. Para 15: merge runs (remove proofErr around "I'd"), no text change
$d.Content.Find.Execute("I’d recommend watching this movie not only once, but multiple times to really get a sense of how the movie conveys emotions and characters.", $true, $false, $false, $false, $false, $true, 1, $false, "I’d recommend watching this movie not only once, but multiple times to really get a sense of how the movie conveys emotions and characters.", 2) | Out-Null

Write-Output "done part 3"
